$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # survey
$ws3 = $wb.Worksheets.Item(3)   # settings

# --- Content fixes ----------------------------------------------------
# The "order number" instructions on the survey sheet had their English
# and French labels swapped; put them back the right way round.
$ws1.Range("C10").Value = "Enter the order number"
$ws1.Range("E10").Value = "Entrer le numéro d'ordre"

# Bump the form title / form id from v2 to v3 on the settings sheet.
$ws3.Range("A2").Value = "(May 2021) 2. Côte d'Ivoire -  Pre TAS FL Formulaire Participants V3"
$ws3.Range("B2").Value = "ci_lf_pretas_2_participant_202105_v3"

# --- View / selection state --------------------------------------------
# "survey" is no longer the selected/active tab; its frozen pane scrolls
# back to column C and the selected cell moves to C25.
$ws1.Activate()
$ws1.Range("C25").Select()

# "settings" becomes the selected/active tab, with B2 selected.
$ws3.Activate()
$ws3.Range("B2").Select()
